$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "总计"
$ws2 = $wb.Worksheets.Item(2)   # currently "2022-Q1" (holds the Q1 fund data)

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q1" sheet so its original data survives as a
#    brand-new sheet placed right after it. That duplicate will keep the name
#    "2022-Q1" (with its original data), while the original sheet object is
#    turned into "2022-Q4" and gets fresh data.
# ---------------------------------------------------------------------------
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)   # duplicate -> will become "2022-Q1"

# Rename: free up "2022-Q1" on the duplicate, relabel the original as "2022-Q4"
$ws2.Name = "2022-Q4"
$ws3.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. Re-style the header row (B1:H1) and the index column (A2:A3) of the new
#    "2022-Q4" sheet to match the "s=2" look used on the "总计" sheet (instead
#    of the "s=1" look that the old "2022-Q1" sheet used).
# ---------------------------------------------------------------------------
$ws1.Range("B1:D1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Overwrite the "2022-Q4" sheet's data rows with the new fund holdings.
#    Columns B and D:G must stay textual (matching the source file's
#    inline-string cells), so force text format before assignment.
# ---------------------------------------------------------------------------
$ws2.Range("B2:B3").NumberFormat = "@"
$ws2.Range("D2:G3").NumberFormat = "@"

$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "007315"
$ws2.Cells.Item(2,3).Value = "汇安嘉盈一年持有期债券A"
$ws2.Cells.Item(2,4).Value = "0.16"
$ws2.Cells.Item(2,5).Value = "24.33"
$ws2.Cells.Item(2,6).Value = "1.37"
$ws2.Cells.Item(2,7).Value = "0.0022"
$ws2.Cells.Item(2,8).Value = 2

$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,2).Value = "010270"
$ws2.Cells.Item(3,3).Value = "汇安嘉盈一年持有期债券C"
$ws2.Cells.Item(3,4).Value = "0.15"
$ws2.Cells.Item(3,5).Value = "24.33"
$ws2.Cells.Item(3,6).Value = "1.37"
$ws2.Cells.Item(3,7).Value = "0.0021"
$ws2.Cells.Item(3,8).Value = 2

# ---------------------------------------------------------------------------
# 4. Update the "总计" overview sheet: insert a new row for "2022-Q4" above
#    the existing "2022-Q1" row (which shifts down to row 3 and gets its
#    index bumped from 0 to 1).
# ---------------------------------------------------------------------------
# Give row 3's index cell (A3) the same formatting as A2 before it is reused.
$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)

# Push the existing "2022-Q1" row down to row 3.
$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "2022-Q1"
$ws1.Cells.Item(3,3).Value = 1
$ws1.Cells.Item(3,4).Value = 0

# Write the new "2022-Q4" row into row 2.
$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q4"
$ws1.Cells.Item(2,3).Value = 2
$ws1.Cells.Item(2,4).Value = 0

# Restore the original active sheet ("总计"), since copying/renaming sheets
# along the way shifts the active tab as a side effect.
$ws1.Activate()
